$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 20) mirroring the existing rows' layout/format.
$row = 20
$prevRow = $row - 1

$ws.Rows.Item($prevRow).Copy()
$ws.Rows.Item($row).PasteSpecial()

$ws.Cells.Item($row, 1).Value = 42625.883356481485
$ws.Cells.Item($row, 2).Value = 24
$ws.Cells.Item($row, 3).Value = 59
$ws.Cells.Item($row, 4).Value = 34
$ws.Cells.Item($row, 5).Value = 92
$ws.Cells.Item($row, 6).Value = 6
$ws.Cells.Item($row, 7).Value = 16326
$ws.Cells.Item($row, 8).Value = 13467
$ws.Cells.Item($row, 9).Value = 2137
$ws.Cells.Item($row, 10).Value = 374
$ws.Cells.Item($row, 11).Value = 215
$ws.Cells.Item($row, 12).Value = 40
$ws.Cells.Item($row, 13).Value = 3
$ws.Cells.Item($row, 14).Value = "Noun"
